# Update cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.227.18"
$ws.Range("E2").Value = "  -4.82%  "
$ws.Range("D3").Value = "3.402.97"
$ws.Range("E3").Value = "  -4.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'559.32"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'172.30"
$ws.Range("E6").Value = "  -9.63%  "
$ws.Range("D7").Value = "'0.616"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.617"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").Value = "'0.152"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").Value = "'56.09"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "'0.0000267"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("D13").Value = "'8.98"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D14").Value = "3.952.70"
$ws.Range("E14").Value = "  -3.82%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.407.26"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.119"
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "'17.91"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "'11.73"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").Value = "64.164.20"
$ws.Range("E19").Value = "  -4.86%  "
$ws.Range("D20").Value = "'0.983"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "'408.02"
$ws.Range("E21").Value = "  -4.81%  "
$ws.Range("D22").Value = "'4.09"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'4.41"
$ws.Range("E23").Value = "  +5.90%  "
$ws.Range("D24").Value = "'13.24"
$ws.Range("E24").Value = "  +7.35%  "
$ws.Range("D25").Value = "'82.57"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").Value = "'10.73"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Value = "'2.74"
$ws.Range("E27").Value = "  -6.01%  "
$ws.Range("D28").Value = "'8.77"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").Value = "'29.49"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("D30").Value = "'6.58"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").Value = "'587.98"
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("D32").Value = "'11.43"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").Value = "'58.89"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("D35").Value = "'0.151"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'36.05"
$ws.Range("E37").Value = "  -6.52%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.370"
$ws.Range("E38").Value = "  -4.89%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "3.183.80"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.39"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0730"
$ws.Range("E41").Value = "  -10.72%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'2.86"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = "  -5.63%  "
$ws.Range("D45").Value = "'3.26"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("D46").Value = "'0.0404"
$ws.Range("E46").Value = "  -3.97%  "
$ws.Range("D47").Value = "'2.61"
$ws.Range("E47").Value = "  -6.15%  "
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'134.25"
$ws.Range("E49").Value = "  -6.47%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'8.22"
$ws.Range("E50").Value = "  -5.07%  "
$ws.Range("E51").Value = "  +2.25%  "
